# Auto-generated edit script: updates 26 data rows across 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 481.47223
$ws.Range("I121").Value = 488.57144
$ws.Range("J121").Value = 479.7586
$ws.Range("K121").Value = 1465.71432
$ws.Range("L121").Value = 1439.2758
$ws.Range("M121").Value = 281.28568
$ws.Range("N121").Value = -4933.275799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 88067.5
$ws.Range("I2").Value = 88067.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 88067.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -87954.5
$ws.Range("N2").ClearContents()
$ws.Range("H25").Value = 500
$ws.Range("I25").Value = 500
$ws.Range("K25").Value = 500
$ws.Range("M25").Value = -98
$ws.Range("H74").Value = 13010210
$ws.Range("I74").Value = 941.52
$ws.Range("J74").Value = 33337192
$ws.Range("K74").Value = 941.52
$ws.Range("L74").Value = 33337192
$ws.Range("M74").Value = -67.51999999999998
$ws.Range("N74").Value = -33338940
$ws.Range("H77").Value = 13010210
$ws.Range("I77").Value = 941.52
$ws.Range("J77").Value = 33337192
$ws.Range("K77").Value = 4707.6
$ws.Range("L77").Value = 166685960
$ws.Range("M77").Value = -339.6000000000004
$ws.Range("N77").Value = -166694696
$ws.Range("H116").Value = 88067.5
$ws.Range("I116").Value = 88067.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 88067.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -85773.5
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 88067.5
$ws.Range("I3").Value = 88067.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 88067.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -87953.5
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 843.7
$ws.Range("I16").Value = 883.2143
$ws.Range("J16").Value = 751.5
$ws.Range("K16").Value = 883.2143
$ws.Range("L16").Value = 751.5
$ws.Range("M16").Value = -596.2143
$ws.Range("N16").Value = -1325.5
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 3207738
$ws.Range("I31").Value = 3624451.8
$ws.Range("J31").Value = 12933.333
$ws.Range("K31").Value = 3624451.8
$ws.Range("L31").Value = 12933.333
$ws.Range("M31").Value = -3624156.8
$ws.Range("N31").Value = -13523.333
$ws.Range("H34").Value = 3207738
$ws.Range("I34").Value = 3624451.8
$ws.Range("J34").Value = 12933.333
$ws.Range("K34").Value = 3624451.8
$ws.Range("L34").Value = 12933.333
$ws.Range("M34").Value = -3624249.8
$ws.Range("N34").Value = -13337.333
$ws.Range("H52").Value = 27620.834
$ws.Range("J52").Value = 27620.834
$ws.Range("L52").Value = 27620.834
$ws.Range("N52").Value = -28208.834
$ws.Range("H113").Value = 843.7
$ws.Range("I113").Value = 883.2143
$ws.Range("J113").Value = 751.5
$ws.Range("K113").Value = 883.2143
$ws.Range("L113").Value = 751.5
$ws.Range("M113").Value = 1286.7857
$ws.Range("N113").Value = -5091.5
$ws.Range("H122").Value = 3279.6511
$ws.Range("I122").Value = 4175.923
$ws.Range("J122").Value = 1908.8823
$ws.Range("K122").Value = 12527.769
$ws.Range("L122").Value = 5726.6469
$ws.Range("M122").Value = -10077.769
$ws.Range("N122").Value = -10626.6469
$ws.Range("H129").Value = 52829.75
$ws.Range("J129").Value = 52829.75
$ws.Range("L129").Value = 52829.75
$ws.Range("N129").Value = -62829.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3032409
$ws.Range("I5").Value = 7693021
$ws.Range("J5").Value = 1737794.6
$ws.Range("K5").Value = 23079063
$ws.Range("L5").Value = 5213383.800000001
$ws.Range("M5").Value = -23078951
$ws.Range("N5").Value = -5213607.800000001
$ws.Range("H46").Value = 2481.3333
$ws.Range("I46").Value = 1111
$ws.Range("J46").Value = 3166.5
$ws.Range("K46").Value = 3333
$ws.Range("L46").Value = 9499.5
$ws.Range("M46").Value = -3242
$ws.Range("N46").Value = -9681.5
$ws.Range("H120").Value = 166667680
$ws.Range("I120").Value = 166667680
$ws.Range("K120").Value = 500003040
$ws.Range("M120").Value = -499998202
$ws.Range("H135").Value = 3032409
$ws.Range("I135").Value = 7693021
$ws.Range("J135").Value = 1737794.6
$ws.Range("K135").Value = 69237189
$ws.Range("L135").Value = 15640151.4
$ws.Range("M135").Value = -69234654
$ws.Range("N135").Value = -15645221.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 32233.334
$ws.Range("J113").Value = 92000
$ws.Range("L113").Value = 92000
$ws.Range("N113").Value = -96340
$ws.Range("H129").Value = 34254
$ws.Range("J129").Value = 34254
$ws.Range("L129").Value = 34254
$ws.Range("N129").Value = -44254

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2343.7222
$ws.Range("I68").Value = 1524
$ws.Range("J68").Value = 3631.8572
$ws.Range("K68").Value = 1524
$ws.Range("L68").Value = 3631.8572
$ws.Range("M68").Value = -775
$ws.Range("N68").Value = -5129.8572
$ws.Range("H71").Value = 2343.7222
$ws.Range("I71").Value = 1524
$ws.Range("J71").Value = 3631.8572
$ws.Range("K71").Value = 7620
$ws.Range("L71").Value = 18159.286
$ws.Range("M71").Value = -3876
$ws.Range("N71").Value = -25647.286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 8680
$ws.Range("J42").Value = 8680
$ws.Range("L42").Value = 8680
$ws.Range("N42").Value = -9436
$ws.Range("H126").Value = 1375
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1375
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 4125
$ws.Range("N126").Value = -9065
$ws.Range("M126").ClearContents()

Write-Host "Edit complete: 26 rows updated across 8 sheets."
